$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "88.331.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.262.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "626.02"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.407"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +12.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.710"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +16.43%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.255.29"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.564"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.68%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +12.03%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.11"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.850.87"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.153.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.287.56"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.12"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.98"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.31"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.31"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.40"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.19%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000142"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +12.55%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.418.88"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Litecoin"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.92"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.181"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -17.11%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.81"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "562.47"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.63%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.36"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -15.69%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.95"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.05"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.140"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.06%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.85"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.82"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.20"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.84%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.399"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.07%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.13%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "149.37"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.16%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "178.50"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.46%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.04"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.17%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.133"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +17.46%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.34"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.72%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.24"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.629"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.57%  "
